$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the site/complex labels:
# B9 used to show "COMPLEX, STA. ROSA, LAGUNA" -> now shows "COMPLEX TECH2"
$ws.Range("B9").Value = "COMPLEX TECH2"
# A10 used to show "COMPLEX TECH" -> now shows "COMPLEX, STA. ROSA, LAGUNA"
$ws.Range("A10").Value = "COMPLEX, STA. ROSA, LAGUNA"

# Update the two date/time serials (keep as raw numeric serials so the
# existing date number-format on these cells is preserved).
$ws.Range("H9").Value = 45317.66666666667
$ws.Range("B16").Value = 45298.25072916667
